$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93 (shifts existing rows 93:233 down to 94:234)
$ws.Rows.Item(93).Insert()

# Populate the new row 93 with a new week's data, mirroring the layout of the
# surrounding rows (same market/category/variety/quality/unit/origin/class),
# with its own date, volume and price figures.
$ws.Range("A93").Value = 10
$ws.Range("B93").Value = "Vega Modelo de Temuco"
$ws.Range("C93").Value = "La Araucanía"
$ws.Range("D93").Value = 44799
$ws.Range("E93").Value = 9
$ws.Range("F93").Value = 100112005
$ws.Range("G93").Value = "Puerro"
$ws.Range("H93").Value = "Azul de Maquehue"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 40
$ws.Range("K93").Value = 17000
$ws.Range("L93").Value = 17000
$ws.Range("M93").Value = 17000
$ws.Range("N93").Value = "$/docena de paquetes"
$ws.Range("O93").Value = "Provincia de Cautín"
$ws.Range("P93").Value = 1417
$ws.Range("Q93").Value = 12
$ws.Range("R93").Value = "Hortaliza"
